# Auto-generated edit script: updates computed profit figures on several
# worksheets' leve-profit data rows, per the scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$edits = @(
    @("H43", "1655.1666"),
    @("J43", "1637.2"),
    @("L43", "1637.2"),
    @("N43", "-1775.2"),
    @("H62", "2761"),
    @("I62", "2726.25"),
    @("K62", "2726.25"),
    @("N62", "-4148"),
    @("M62", "-2102.25"),
    @("J62", "2900"),
    @("L62", "2900"),
    @("M65", "-10511.25"),
    @("L65", "14500"),
    @("K65", "13631.25"),
    @("I65", "2726.25"),
    @("J65", "2900"),
    @("H65", "2761"),
    @("N65", "-20740"),
    @("K101", "1411.00002"),
    @("I101", "470.33334"),
    @("J101", "0"),
    @("M101", "210.9999800000001"),
    @("N101", $null),
    @("H101", "470.33334"),
    @("L101", "0"),
    @("N116", "-11800.9165"),
    @("H116", "4059"),
    @("L116", "4916.9165"),
    @("J116", "4916.9165"),
    @("N121", "-8999.9231"),
    @("J121", "1835.3077"),
    @("H121", "1683.9333"),
    @("L121", "5505.9231"),
    @("H129", "855.4375"),
    @("J129", "855.4375"),
    @("N129", "-12566.3125"),
    @("L129", "2566.3125"),
    @("I132", "4335.2354"),
    @("H132", "4295.1904"),
    @("K132", "13005.7062"),
    @("N132", "-17435"),
    @("J132", "4125"),
    @("M132", "-10475.7062"),
    @("L132", "12375"),
    @("J137", "168837.17"),
    @("M137", "-2550.857400000001"),
    @("L137", "506511.51"),
    @("N137", "-511611.51"),
    @("I137", "1700.2858"),
    @("K137", "5100.857400000001"),
    @("H137", "78840.38"),
    @("I138", "574.6"),
    @("M138", "3416.2"),
    @("J138", "3383.6667"),
    @("L138", "10151.0001"),
    @("N138", "-20431.0001"),
    @("H138", "2013.3903"),
    @("K138", "1723.8"),
    @("H141", "3031.6667"),
    @("M141", "990.0002000000004"),
    @("K141", "4189.9998"),
    @("I141", "1396.6666")
)
foreach ($edit in $edits) {
    $cellRef = $edit[0]
    $newVal = $edit[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = [double]$newVal
    }
}

$ws = $wb.Worksheets.Item("ARM")
$edits = @(
    @("H32", "2372.9412"),
    @("M32", "-1705.0834"),
    @("I32", "1992.0834"),
    @("K32", "1992.0834"),
    @("I61", "727.63416"),
    @("L61", "3636.6667"),
    @("M61", "-515.63416"),
    @("K61", "727.63416"),
    @("J61", "3636.6667"),
    @("N61", "-4060.6667"),
    @("H61", "1506.8392"),
    @("J136", "3636.6667"),
    @("M136", "367.0975200000003"),
    @("I136", "727.63416"),
    @("N136", "-16010.0001"),
    @("H136", "1506.8392"),
    @("K136", "2182.90248"),
    @("L136", "10910.0001")
)
foreach ($edit in $edits) {
    $cellRef = $edit[0]
    $newVal = $edit[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = [double]$newVal
    }
}

$ws = $wb.Worksheets.Item("BSM")
$edits = @(
    @("I94", "1799.7778"),
    @("K94", "1799.7778"),
    @("H94", "3273.4666"),
    @("M94", "-1348.7778")
)
foreach ($edit in $edits) {
    $cellRef = $edit[0]
    $newVal = $edit[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = [double]$newVal
    }
}

$ws = $wb.Worksheets.Item("CRP")
$edits = @(
    @("L31", "3964"),
    @("H31", "2799.5151"),
    @("N31", "-4554"),
    @("I31", "2134.0952"),
    @("M31", "-1839.0952"),
    @("J31", "3964"),
    @("K31", "2134.0952"),
    @("J34", "3964"),
    @("K34", "2134.0952"),
    @("N34", "-4368"),
    @("L34", "3964"),
    @("M34", "-1932.0952"),
    @("H34", "2799.5151"),
    @("I34", "2134.0952"),
    @("H122", "988.7778"),
    @("I122", "988.7778"),
    @("M122", "-516.3334"),
    @("K122", "2966.3334"),
    @("I132", "2087.2"),
    @("H132", "3198.1667"),
    @("K132", "6261.599999999999"),
    @("M132", "-3731.599999999999"),
    @("K134", "2652"),
    @("J134", "1613"),
    @("M134", "-117"),
    @("H134", "1182.2273"),
    @("N134", "-9909"),
    @("L134", "4839"),
    @("I134", "884")
)
foreach ($edit in $edits) {
    $cellRef = $edit[0]
    $newVal = $edit[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = [double]$newVal
    }
}

$ws = $wb.Worksheets.Item("CUL")
$edits = @(
    @("N5", "-6599"),
    @("J5", "2125"),
    @("L5", "6375"),
    @("H5", "1365.0588"),
    @("I68", "0"),
    @("K68", "0"),
    @("H68", "33934.332"),
    @("L68", "101802.996"),
    @("M68", $null),
    @("J68", "33934.332"),
    @("N68", "-103424.996"),
    @("I71", "0"),
    @("H71", "33934.332"),
    @("N71", "-313520.988"),
    @("J71", "33934.332"),
    @("L71", "305408.988"),
    @("M71", $null),
    @("K71", "0"),
    @("K131", "1180.00002"),
    @("J131", "785.3936"),
    @("L131", "2356.1808"),
    @("H131", "761.87"),
    @("I131", "393.33334"),
    @("M131", "3859.99998"),
    @("N131", "-12436.1808"),
    @("H135", "1365.0588"),
    @("J135", "2125"),
    @("N135", "-24195"),
    @("L135", "19125")
)
foreach ($edit in $edits) {
    $cellRef = $edit[0]
    $newVal = $edit[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = [double]$newVal
    }
}

$ws = $wb.Worksheets.Item("LTW")
$edits = @(
    @("I68", "2054.2856"),
    @("K68", "2054.2856"),
    @("H68", "4103.9414"),
    @("M68", "-1305.2856"),
    @("I71", "2054.2856"),
    @("H71", "4103.9414"),
    @("M71", "-6527.428"),
    @("K71", "10271.428"),
    @("L122", "16608"),
    @("H122", "1511969.2"),
    @("N122", "-21508"),
    @("I122", "2181495.2"),
    @("J122", "5536"),
    @("M122", "-6542035.600000001"),
    @("K122", "6544485.600000001"),
    @("I132", "2030.6666"),
    @("H132", "2835.182"),
    @("K132", "6091.9998"),
    @("M132", "-3561.9998")
)
foreach ($edit in $edits) {
    $cellRef = $edit[0]
    $newVal = $edit[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = [double]$newVal
    }
}

$ws = $wb.Worksheets.Item("WVR")
$edits = @(
    @("H107", "6495470.5"),
    @("K107", "4500"),
    @("N107", "-22737237"),
    @("M107", "-2580"),
    @("L107", "22733397"),
    @("J107", "7577799"),
    @("I107", "1500"),
    @("I126", "588.1875"),
    @("K126", "1764.5625"),
    @("M126", "705.4375"),
    @("H126", "1162.4286")
)
foreach ($edit in $edits) {
    $cellRef = $edit[0]
    $newVal = $edit[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = [double]$newVal
    }
}
